$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ALMACENES"
$ws.Range("B2").Value = "CLIENTE B1"
$ws.Range("C2").Value = "00000015"
$ws.Range("D2").Value = "- - -3654789"

$ws.Range("D2").Select()
